$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove cells that became empty (naive forecaster bug fix)
$cellsToClear = @("C2", "E2", "C3", "C4")
foreach ($cellRef in $cellsToClear) {
    $ws.Range($cellRef).ClearContents()
}

# Update recalculated forecast values (floating point precision refresh)
$ws.Range("E3").Value = 8.045645122021927
$ws.Range("E4").Value = 6.772115316529947
$ws.Range("C5").Value = -5.232639093663827
$ws.Range("C6").Value = -7.266312015249799
$ws.Range("C7").Value = 2.124540184802992
$ws.Range("C9").Value = 8.081020954067753
$ws.Range("C10").Value = 9.469137444079955
$ws.Range("C11").Value = 4.489210662380971
$ws.Range("C12").Value = 3.358206407534969
$ws.Range("E12").Value = 5.745831525574463
$ws.Range("C13").Value = -0.8752093743685241
$ws.Range("C16").Value = 3.901355411819685
$ws.Range("E16").Value = 3.690459963535031
$ws.Range("C17").Value = 4.818339085077561
$ws.Range("C21").Value = 4.613634856640769
$ws.Range("C22").Value = 5.246209615995689
$ws.Range("C23").Value = 6.155351106582851
$ws.Range("C26").Value = 4.862559663742938
$ws.Range("E26").Value = 4.636196713604379
$ws.Range("E28").Value = 4.838485897465628
$ws.Range("C30").Value = 2.76474001115945
$ws.Range("E31").Value = 2.845541644111549
$ws.Range("C33").Value = -9.509392583043464
$ws.Range("C34").Value = -7.260793671746447
$ws.Range("E35").Value = -1.005130103122098
$ws.Range("E36").Value = 0.869978169785246
$ws.Range("C38").Value = 4.097586525396246
$ws.Range("C39").Value = 8.215174201986319
$ws.Range("E39").Value = 2.364846754700167
$ws.Range("E40").Value = 2.551560717335266
$ws.Range("E41").Value = 3.242110390729347
$ws.Range("C42").Value = 7.824284864703768
$ws.Range("C46").Value = -1.245022353133318
$ws.Range("C47").Value = -3.067646799613699
$ws.Range("E47").Value = 0.4897355149953819
$ws.Range("C48").Value = -2.447533648174649
$ws.Range("C49").Value = -2.520879465820702
$ws.Range("E49").Value = -0.01136840354140078
$ws.Range("E51").Value = -0.3110100908356728
$ws.Range("C52").Value = 1.038949519463617
$ws.Range("E52").Value = -1.220869074712128
